$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need NumberFormat forced
# to Text ("@") before assignment, otherwise Excel auto-converts them to
# numeric values and things like trailing zeros / multi-dot "thousands"
# notation (e.g. "23.421.24", "1.000") would be silently lost.

$ws.Range('D2').Value = '23.421.24'
$ws.Range('E2').Value = '  +1.04%  '
$ws.Range('D3').Value = '1.638.37'
$ws.Range('E3').Value = '  +2.31%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('E5').Value = '  +0.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '304.95'
$ws.Range('E6').Value = '  +0.64%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3741'
$ws.Range('E7').Value = '  -1.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '51.92'
$ws.Range('E8').Value = '  -0.20%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3621'
$ws.Range('E9').Value = '  +0.14%  '
$ws.Range('E10').Value = '  -0.79%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08124'
$ws.Range('E11').Value = '  +0.10%  '
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.85'
$ws.Range('E13').Value = '  +0.37%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.602'
$ws.Range('E14').Value = '  +0.28%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001268'
$ws.Range('E15').Value = '  +2.00%  '
$ws.Range('E16').Value = '  -1.59%  '
$ws.Range('D17').Value = '1.637.75'
$ws.Range('E17').Value = '  +2.43%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '94.25'
$ws.Range('E18').Value = '  +0.36%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06905'
$ws.Range('E19').Value = '  +0.28%  '
$ws.Range('E20').Value = '  +0.39%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.503'
$ws.Range('E21').Value = '  -0.49%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.000'
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('D23').Value = '23.434.17'
$ws.Range('E23').Value = '  +1.12%  '
$ws.Range('E24').Value = '  -1.78%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.420'
$ws.Range('E25').Value = '  +1.00%  '
$ws.Range('B26').Value = 'LidoDAOToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.031'
$ws.Range('E26').Value = '  +1.38%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '21.20'
$ws.Range('E27').Value = '  -0.05%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '151.42'
$ws.Range('E28').Value = '  +1.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.317'
$ws.Range('E29').Value = '  +1.35%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '135.79'
$ws.Range('E30').Value = '  +1.58%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.272'
$ws.Range('E31').Value = '  -3.93%  '
$ws.Range('D32').Value = '1.818.30'
$ws.Range('E32').Value = '  +2.25%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.734'
$ws.Range('E33').Value = '  -0.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9551'
$ws.Range('E34').Value = '  -1.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02817'
$ws.Range('E35').Value = '  +3.60%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '10.30'
$ws.Range('E36').Value = '  +0.77%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.07249'
$ws.Range('E37').Value = '  -2.91%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2517'
$ws.Range('E38').Value = '  +0.13%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.08784'
$ws.Range('E39').Value = '  -0.14%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.068'
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('E41').Value = '  +1.21%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.7038'
$ws.Range('E42').Value = '  -0.81%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '12.43'
$ws.Range('E43').Value = '  -0.32%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.01'
$ws.Range('E44').Value = '  +2.49%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6501'
$ws.Range('E45').Value = '  -0.44%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.325'
$ws.Range('E46').Value = '  +0.64%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.000'
$ws.Range('E47').Value = '  +0.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.007'
$ws.Range('E48').Value = '  -0.28%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.07978'
$ws.Range('E49').Value = '  +0.38%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '128.26'
$ws.Range('E50').Value = '  -2.82%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.201'
$ws.Range('E51').Value = '  -0.03%  '
